$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds numeric-looking strings (e.g. "27.297.96",
# "0.9991") that must stay literal text, like the source inline-string
# cells. Mark the column as Text first so Excel does not coerce these
# into numbers/dates when the values are assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.297.96"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3
$ws.Range("D3").Value = "1.826.27"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "313.58"
$ws.Range("E5").Value = "  +0.66%  "

# Row 6
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("D7").Value = "0.4483"
$ws.Range("E7").Value = "  +5.36%  "

# Row 8
$ws.Range("D8").Value = "0.3793"
$ws.Range("E8").Value = "  +3.32%  "

# Row 9
$ws.Range("D9").Value = "0.07513"
$ws.Range("E9").Value = "  +3.80%  "

# Row 10
$ws.Range("D10").Value = "0.8841"
$ws.Range("E10").Value = "  +4.69%  "

# Row 11
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").Value = "1.812.70"
$ws.Range("E12").Value = "  -0.76%  "

# Row 13
$ws.Range("D13").Value = "6.771"
$ws.Range("E13").Value = "  +1.45%  "

# Row 14
$ws.Range("D14").Value = "94.73"
$ws.Range("E14").Value = "  +5.64%  "

# Row 15
$ws.Range("D15").Value = "5.402"
$ws.Range("E15").Value = "  +2.03%  "

# Row 16
$ws.Range("D16").Value = "0.07132"
$ws.Range("E16").Value = "  +1.35%  "

# Row 17
$ws.Range("D17").Value = "0.9984"
$ws.Range("E17").Value = "  -0.39%  "

# Row 18
$ws.Range("D18").Value = "0.000008799"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19
$ws.Range("D19").Value = "0.9989"
$ws.Range("E19").Value = "  -0.16%  "

# Row 20
$ws.Range("D20").Value = "15.16"
$ws.Range("E20").Value = "  +1.82%  "

# Row 21
$ws.Range("D21").Value = "27.335.80"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("D22").Value = "5.257"
$ws.Range("E22").Value = "  +2.41%  "

# Row 23
$ws.Range("D23").Value = "11.00"
$ws.Range("E23").Value = "  +1.74%  "

# Row 24
$ws.Range("D24").Value = "2.046.83"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("E25").Value = "  +0.36%  "

# Row 26
$ws.Range("D26").Value = "2.454"
$ws.Range("E26").Value = "  +8.50%  "

# Row 27
$ws.Range("D27").Value = "152.09"
$ws.Range("E27").Value = "  +0.42%  "

# Row 28
$ws.Range("D28").Value = "18.67"
$ws.Range("E28").Value = "  +2.55%  "

# Row 29
$ws.Range("D29").Value = "5.382"
$ws.Range("E29").Value = "  +2.63%  "

# Row 30
$ws.Range("D30").Value = "118.46"
$ws.Range("E30").Value = "  +1.44%  "

# Row 31
$ws.Range("D31").Value = "0.08859"
$ws.Range("E31").Value = "  +1.32%  "

# Row 32
$ws.Range("D32").Value = "0.7729"
$ws.Range("E32").Value = "  +4.81%  "

# Row 33
$ws.Range("E33").Value = "  +0.71%  "

# Row 34
$ws.Range("E34").Value = "  +3.45%  "

# Row 35
$ws.Range("D35").Value = "2.884"
$ws.Range("E35").Value = "  -0.59%  "

# Row 36
$ws.Range("D36").Value = "0.9978"
$ws.Range("E36").Value = "  -0.26%  "

# Row 37
$ws.Range("D37").Value = "1.111"
$ws.Range("E37").Value = "  +1.42%  "

# Row 38
$ws.Range("D38").Value = "0.01993"
$ws.Range("E38").Value = "  +2.63%  "

# Row 39
$ws.Range("D39").Value = "0.05316"
$ws.Range("E39").Value = "  +1.63%  "

# Row 40
$ws.Range("D40").Value = "7.427"
$ws.Range("E40").Value = "  +1.60%  "

# Row 41
$ws.Range("D41").Value = "0.5349"
$ws.Range("E41").Value = "  +4.75%  "

# Row 42
$ws.Range("D42").Value = "0.1729"
$ws.Range("E42").Value = "  +2.48%  "

# Row 43
$ws.Range("D43").Value = "2.858"
$ws.Range("E43").Value = "  -0.46%  "

# Row 44
$ws.Range("E44").Value = "  +14.10%  "

# Row 45
$ws.Range("D45").Value = "8.799"
$ws.Range("E45").Value = "  +2.67%  "

# Row 46
$ws.Range("D46").Value = "0.5103"
$ws.Range("E46").Value = "  +7.80%  "

# Row 47
$ws.Range("D47").Value = "10.74"
$ws.Range("E47").Value = "  +2.03%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "106.68"
$ws.Range("E48").Value = "  +1.07%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.706"
$ws.Range("E49").Value = "  +3.13%  "

# Row 50
$ws.Range("D50").Value = "0.9978"
$ws.Range("E50").Value = "  -0.27%  "

# Row 51
$ws.Range("D51").Value = "0.06381"
$ws.Range("E51").Value = "  +0.88%  "
